# Applies the cryptos list refresh described in the commit:
# "Updated cryptos list on Tue Jun 20 02:59:52 UTC 2023 with GitHub Actions"
#
# Every data cell in columns B:E is stored as text (string), including values
# that look numeric ("0.9976", "3.943", ...). A plain COM .Value assignment
# auto-detects such strings as numbers, so each write is forced to text by
# prefixing a leading apostrophe (the classic "text-literal" marker) and then
# resetting the cell Style back to "Normal" so no stray number-format/
# quote-prefix style sticks around afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.Value = "'26.844.91"
$cell.Style = "Normal"

# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.Value = "'1.731.64"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 5)
$cell.Value = "'  +0.32%  "
$cell.Style = "Normal"

# Row 4
$cell = $ws.Cells.Item(4, 4)
$cell.Value = "'0.9976"
$cell.Style = "Normal"

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.Value = "'242.27"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 5)
$cell.Value = "'  -0.67%  "
$cell.Style = "Normal"

# Row 6
$cell = $ws.Cells.Item(6, 5)
$cell.Value = "'  -0.17%  "
$cell.Style = "Normal"

# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.Value = "'0.4913"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(7, 5)
$cell.Value = "'  +0.04%  "
$cell.Style = "Normal"

# Row 8
$cell = $ws.Cells.Item(8, 4)
$cell.Value = "'0.2610"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(8, 5)
$cell.Value = "'  -0.27%  "
$cell.Style = "Normal"

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.Value = "'0.06223"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(9, 5)
$cell.Value = "'  +0.24%  "
$cell.Style = "Normal"

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.Value = "'1.737.10"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(10, 5)
$cell.Value = "'  +0.69%  "
$cell.Style = "Normal"

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.Value = "'16.08"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 5)
$cell.Value = "'  +3.44%  "
$cell.Style = "Normal"

# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.Value = "'0.06913"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 5)
$cell.Value = "'  -1.62%  "
$cell.Style = "Normal"

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.Value = "'0.6120"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(13, 5)
$cell.Value = "'  +1.55%  "
$cell.Style = "Normal"

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.Value = "'4.513"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 5)
$cell.Value = "'  -1.29%  "
$cell.Style = "Normal"

# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.Value = "'77.37"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 5)
$cell.Value = "'  -0.01%  "
$cell.Style = "Normal"

# Row 16
$cell = $ws.Cells.Item(16, 5)
$cell.Value = "'  -0.10%  "
$cell.Style = "Normal"

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.Value = "'26.624.54"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(17, 5)
$cell.Value = "'  +0.70%  "
$cell.Style = "Normal"

# Row 18
$cell = $ws.Cells.Item(18, 4)
$cell.Value = "'0.9976"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 5)
$cell.Value = "'  -0.19%  "
$cell.Style = "Normal"

# Row 19
$cell = $ws.Cells.Item(19, 4)
$cell.Value = "'0.000007191"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 5)
$cell.Value = "'  -0.07%  "
$cell.Style = "Normal"

# Row 20
$cell = $ws.Cells.Item(20, 5)
$cell.Value = "'  +0.96%  "
$cell.Style = "Normal"

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.Value = "'1.960.64"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(21, 5)
$cell.Value = "'  +0.79%  "
$cell.Style = "Normal"

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.Value = "'8.578"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(23, 5)
$cell.Value = "'  -0.33%  "
$cell.Style = "Normal"

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.Value = "'5.131"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 5)
$cell.Value = "'  -0.83%  "
$cell.Style = "Normal"

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.Value = "'138.88"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(25, 5)
$cell.Value = "'  +0.87%  "
$cell.Style = "Normal"

# Row 26
$cell = $ws.Cells.Item(26, 4)
$cell.Value = "'15.32"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 5)
$cell.Value = "'  +0.48%  "
$cell.Style = "Normal"

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.Value = "'1.788"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 5)
$cell.Value = "'  +4.91%  "
$cell.Style = "Normal"

# Row 28
$cell = $ws.Cells.Item(28, 2)
$cell.Value = "'Toncoin"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 3)
$cell.Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 4)
$cell.Value = "'1.380"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 5)
$cell.Value = "'  -0.70%  "
$cell.Style = "Normal"

# Row 29
$cell = $ws.Cells.Item(29, 2)
$cell.Value = "'BitcoinCash"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 3)
$cell.Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 4)
$cell.Value = "'106.13"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 5)
$cell.Value = "'  -0.94%  "
$cell.Style = "Normal"

# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.Value = "'3.943"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(30, 5)
$cell.Value = "'  -0.60%  "
$cell.Style = "Normal"

# Row 31
$cell = $ws.Cells.Item(31, 4)
$cell.Value = "'0.08001"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(31, 5)
$cell.Value = "'  +0.33%  "
$cell.Style = "Normal"

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.Value = "'3.686"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(32, 5)
$cell.Value = "'  +0.05%  "
$cell.Style = "Normal"

# Row 34
$cell = $ws.Cells.Item(34, 2)
$cell.Value = "'Frax"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 3)
$cell.Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 4)
$cell.Value = "'0.9974"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 5)
$cell.Value = "'  -0.16%  "
$cell.Style = "Normal"

# Row 35
$cell = $ws.Cells.Item(35, 2)
$cell.Value = "'HuobiToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 3)
$cell.Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 4)
$cell.Value = "'2.603"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 5)
$cell.Value = "'  +0.11%  "
$cell.Style = "Normal"

# Row 36
$cell = $ws.Cells.Item(36, 2)
$cell.Value = "'ARBITRUM"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 3)
$cell.Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 4)
$cell.Value = "'1.010"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(36, 5)
$cell.Value = "'  +0.98%  "
$cell.Style = "Normal"

# Row 37
$cell = $ws.Cells.Item(37, 2)
$cell.Value = "'ImmutableX"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 3)
$cell.Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 4)
$cell.Value = "'0.6252"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(37, 5)
$cell.Value = "'  -0.40%  "
$cell.Style = "Normal"

# Row 38
$cell = $ws.Cells.Item(38, 2)
$cell.Value = "'TrustWalletToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 3)
$cell.Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 4)
$cell.Value = "'0.9395"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(38, 5)
$cell.Value = "'  +3.00%  "
$cell.Style = "Normal"

# Row 39
$cell = $ws.Cells.Item(39, 2)
$cell.Value = "'RenderToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 3)
$cell.Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 4)
$cell.Value = "'2.052"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(39, 5)
$cell.Value = "'  +4.33%  "
$cell.Style = "Normal"

# Row 40
$cell = $ws.Cells.Item(40, 2)
$cell.Value = "'MXToken"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 3)
$cell.Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 4)
$cell.Value = "'2.448"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(40, 5)
$cell.Value = "'  +2.36%  "
$cell.Style = "Normal"

# Row 41
$cell = $ws.Cells.Item(41, 2)
$cell.Value = "'PaxDollar"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 3)
$cell.Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 4)
$cell.Value = "'0.9983"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 5)
$cell.Value = "'  -0.17%  "
$cell.Style = "Normal"

# Row 42
$cell = $ws.Cells.Item(42, 2)
$cell.Value = "'VeChain"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 3)
$cell.Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 4)
$cell.Value = "'0.01506"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 5)
$cell.Value = "'  +1.38%  "
$cell.Style = "Normal"

# Row 43
$cell = $ws.Cells.Item(43, 2)
$cell.Value = "'FraxShare"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 3)
$cell.Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 4)
$cell.Value = "'5.655"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 5)
$cell.Value = "'  +3.88%  "
$cell.Style = "Normal"

# Row 44
$cell = $ws.Cells.Item(44, 2)
$cell.Value = "'Quant"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 3)
$cell.Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 4)
$cell.Value = "'99.67"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(44, 5)
$cell.Value = "'  -0.19%  "
$cell.Style = "Normal"

# Row 45
$cell = $ws.Cells.Item(45, 2)
$cell.Value = "'TheSandbox"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 3)
$cell.Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 4)
$cell.Value = "'0.3869"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(45, 5)
$cell.Value = "'  +0.37%  "
$cell.Style = "Normal"

# Row 46
$cell = $ws.Cells.Item(46, 2)
$cell.Value = "'Aptos"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 3)
$cell.Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 4)
$cell.Value = "'6.949"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(46, 5)
$cell.Value = "'  +3.17%  "
$cell.Style = "Normal"

# Row 47
$cell = $ws.Cells.Item(47, 2)
$cell.Value = "'Algorand"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 3)
$cell.Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 4)
$cell.Value = "'0.1162"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(47, 5)
$cell.Value = "'  +0.40%  "
$cell.Style = "Normal"

# Row 48
$cell = $ws.Cells.Item(48, 2)
$cell.Value = "'Cronos"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 3)
$cell.Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 4)
$cell.Value = "'0.05388"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 5)
$cell.Value = "'  +0.39%  "
$cell.Style = "Normal"

# Row 49
$cell = $ws.Cells.Item(49, 2)
$cell.Value = "'EnergySwap"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 3)
$cell.Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 4)
$cell.Value = "'7.966"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 5)
$cell.Value = "'  +3.59%  "
$cell.Style = "Normal"

# Row 50
$cell = $ws.Cells.Item(50, 2)
$cell.Value = "'Elrond"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 3)
$cell.Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 4)
$cell.Value = "'30.26"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 5)
$cell.Value = "'  +0.39%  "
$cell.Style = "Normal"

# Row 51
$cell = $ws.Cells.Item(51, 2)
$cell.Value = "'NEARProtocol"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 3)
$cell.Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(51, 4)
$cell.Value = "'1.245"
$cell.Style = "Normal"

